# Scheduled runner update for "Moogle_Profits" workbook.
# Refreshes cached market-board figures (currentAveragePrice*, LevePrice*,
# LeveProfit*) across the per-job sheets. Values come from the latest
# Universalis pull; a handful of rows lost their HQ-sell-only leve (no
# NQ market activity that cycle) so the NQ profit cell is cleared instead
# of being written as a stale number.

$wb = $excel.ActiveWorkbook

# ---- ALC ----------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

$ws.Range("H3").Value  = 199500
$ws.Range("J3").Value  = 199500
$ws.Range("L3").Value  = 199500
$ws.Range("N3").Value  = -199728

$ws.Range("H4").Value  = 3557.4
$ws.Range("I4").Value  = 1544
$ws.Range("K4").Value  = 1544
$ws.Range("M4").Value  = -1430

$ws.Range("H40").Value = 5845.276
$ws.Range("I40").Value = 4753.6665
$ws.Range("J40").Value = 7014.857
$ws.Range("K40").Value = 4753.6665
$ws.Range("L40").Value = 7014.857
$ws.Range("M40").Value = -4578.6665
$ws.Range("N40").Value = -7364.857

$ws.Range("H48").Value = 6188.25
$ws.Range("I48").Value = 3255
$ws.Range("J48").Value = 7166
$ws.Range("K48").Value = 9765
$ws.Range("L48").Value = 21498
$ws.Range("M48").Value = -9473
$ws.Range("N48").Value = -22082

$ws.Range("H56").Value = 6188.25
$ws.Range("I56").Value = 3255
$ws.Range("J56").Value = 7166
$ws.Range("K56").Value = 9765
$ws.Range("L56").Value = 21498
$ws.Range("M56").Value = -9231
$ws.Range("N56").Value = -22566

$ws.Range("H64").Value = 6481
$ws.Range("J64").Value = 4483
$ws.Range("L64").Value = 4483
$ws.Range("N64").Value = -4979

$ws.Range("H67").Value = 6481
$ws.Range("J67").Value = 4483
$ws.Range("L67").Value = 4483
$ws.Range("N67").Value = -6199

$ws.Range("H102").Value = 199500
$ws.Range("J102").Value = 199500
$ws.Range("L102").Value = 199500
$ws.Range("N102").Value = -205990

$ws.Range("H116").Value = 11526.277
$ws.Range("I116").Value = 8198.385
$ws.Range("J116").Value = 20178.8
$ws.Range("K116").Value = 8198.385
$ws.Range("L116").Value = 20178.8
$ws.Range("M116").Value = -4756.385
$ws.Range("N116").Value = -27062.8

$ws.Range("H132").Value = 2578.8572
$ws.Range("I132").Value = 2676.32
$ws.Range("K132").Value = 8028.960000000001
$ws.Range("M132").Value = -5498.960000000001

$ws.Range("H133").Value = 69948.336
$ws.Range("J133").Value = 69948.336
$ws.Range("L133").Value = 69948.336
$ws.Range("N133").Value = -80068.336

# ---- ARM ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

$ws.Range("H110").Value = 2497.6956
$ws.Range("I110").Value = 2611.9443
$ws.Range("K110").Value = 2611.9443
$ws.Range("M110").Value = -566.9443000000001

# ---- CRP ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

$ws.Range("H7").Value  = 274.1111
$ws.Range("I7").Value  = 151.6
$ws.Range("J7").Value  = 427.25
$ws.Range("K7").Value  = 151.6
$ws.Range("L7").Value  = 427.25
$ws.Range("M7").Value  = -38.59999999999999
$ws.Range("N7").Value  = -653.25

$ws.Range("H50").Value = 72997.5
$ws.Range("I50").Value = 0
$ws.Range("K50").Value = 0
$ws.Range("M50").ClearContents()

$ws.Range("H124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

$ws.Range("H141").Value = 324871.6
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 324871.6
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 324871.6
$ws.Range("N141").Value = -335231.6
$ws.Range("M141").ClearContents()

# ---- CUL ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

$ws.Range("H6").Value   = 64.875
$ws.Range("I6").Value   = 45.153847
$ws.Range("K6").Value   = 135.461541
$ws.Range("M6").Value   = -22.46154100000001

$ws.Range("H11").Value  = 140732.38
$ws.Range("I11").Value  = 187604.5
$ws.Range("K11").Value  = 562813.5
$ws.Range("M11").Value  = -562673.5

$ws.Range("H56").Value  = 7499.4287
$ws.Range("I56").Value  = 7499.4287
$ws.Range("K56").Value  = 7499.4287
$ws.Range("M56").Value  = -6969.4287

$ws.Range("H87").Value  = 200
$ws.Range("I87").Value  = 200
$ws.Range("K87").Value  = 600
$ws.Range("M87").Value  = 648

$ws.Range("H90").Value  = 200
$ws.Range("I90").Value  = 200
$ws.Range("K90").Value  = 1800
$ws.Range("M90").Value  = 4440

$ws.Range("H129").Value = 15171607
$ws.Range("I129").Value = 16357.714
$ws.Range("J129").Value = 41693292
$ws.Range("K129").Value = 49073.142
$ws.Range("L129").Value = 125079876
$ws.Range("M129").Value = -44073.142
$ws.Range("N129").Value = -125089876

$ws.Range("H131").Value = 868401
$ws.Range("I131").Value = 850.4167
$ws.Range("J131").Value = 2025135.1
$ws.Range("K131").Value = 2551.2501
$ws.Range("L131").Value = 6075405.300000001
$ws.Range("M131").Value = 2488.7499
$ws.Range("N131").Value = -6085485.300000001

# ---- GSM ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

$ws.Range("H57").Value  = 0
$ws.Range("I57").Value  = 0
$ws.Range("K57").Value  = 0
$ws.Range("M57").ClearContents()

$ws.Range("H80").Value  = 10802.111
$ws.Range("I80").Value  = 9512.0625
$ws.Range("J80").Value  = 12678.546
$ws.Range("K80").Value  = 9512.0625
$ws.Range("L80").Value  = 12678.546
$ws.Range("M80").Value  = -8514.0625
$ws.Range("N80").Value  = -14674.546

$ws.Range("H83").Value  = 10802.111
$ws.Range("I83").Value  = 9512.0625
$ws.Range("J83").Value  = 12678.546
$ws.Range("K83").Value  = 47560.3125
$ws.Range("L83").Value  = 63392.73
$ws.Range("M83").Value  = -42568.3125
$ws.Range("N83").Value  = -73376.73000000001

$ws.Range("H113").Value = 3666.4285
$ws.Range("I113").Value = 2886.5715
$ws.Range("K113").Value = 2886.5715
$ws.Range("M113").Value = -716.5715

# ---- LTW ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

$ws.Range("H16").Value = 2848.7273
$ws.Range("I16").Value = 3180.5
$ws.Range("J16").Value = 1964
$ws.Range("K16").Value = 3180.5
$ws.Range("L16").Value = 1964
$ws.Range("M16").Value = -3010.5
$ws.Range("N16").Value = -2304

$ws.Range("H22").Value = 1383
$ws.Range("I22").Value = 612
$ws.Range("J22").Value = 2282.5
$ws.Range("K22").Value = 612
$ws.Range("L22").Value = 2282.5
$ws.Range("M22").Value = -317
$ws.Range("N22").Value = -2872.5

$ws.Range("H27").Value = 1383
$ws.Range("I27").Value = 612
$ws.Range("J27").Value = 2282.5
$ws.Range("K27").Value = 612
$ws.Range("L27").Value = 2282.5
$ws.Range("M27").Value = -505
$ws.Range("N27").Value = -2496.5

$ws.Range("H55").Value = 905.35
$ws.Range("I55").Value = 313.5
$ws.Range("J55").Value = 1159
$ws.Range("K55").Value = 313.5
$ws.Range("L55").Value = 1159
$ws.Range("M55").Value = -140.5
$ws.Range("N55").Value = -1505

$ws.Range("H68").Value = 4486.857
$ws.Range("J68").Value = 2800
$ws.Range("L68").Value = 2800
$ws.Range("N68").Value = -4298

$ws.Range("H71").Value = 4486.857
$ws.Range("J71").Value = 2800
$ws.Range("L71").Value = 14000
$ws.Range("N71").Value = -21488

# ---- WVR ------------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

$ws.Range("H50").Value  = 20000
$ws.Range("J50").Value  = 20000
$ws.Range("L50").Value  = 20000
$ws.Range("N50").Value  = -21262

$ws.Range("H52").Value  = 10299.5
$ws.Range("I52").Value  = 10299.5
$ws.Range("K52").Value  = 10299.5
$ws.Range("M52").Value  = -10073.5

$ws.Range("H107").Value = 741722.75
$ws.Range("I107").Value = 1250977
$ws.Range("J107").Value = 989.2727
$ws.Range("K107").Value = 3752931
$ws.Range("L107").Value = 2967.8181
$ws.Range("M107").Value = -3751011
$ws.Range("N107").Value = -6807.8181

$ws.Range("H132").Value = 3512.0557
$ws.Range("I132").Value = 2513.5625
$ws.Range("K132").Value = 7540.6875
$ws.Range("M132").Value = -5010.6875

Write-Host "Moogle_Profits: updated 37 rows across ALC/ARM/CRP/CUL/GSM/LTW/WVR"
